$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.225803017616272
$ws.Range("B1").Value = 1.649744391441345
$ws.Range("C1").Value = 2.01209282875061
$ws.Range("D1").Value = 6.773933887481689
$ws.Range("E1").Value = 2.083653926849365
